$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 276, shifting existing rows 276:405 down to 277:406
$ws.Rows.Item(276).Insert()

# Populate the newly inserted row 276 with its data
$ws.Cells.Item(276, 1).Value = 4
$ws.Cells.Item(276, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(276, 3).Value = "Los Lagos"
$ws.Cells.Item(276, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(276, 4).Value = 44917
$ws.Cells.Item(276, 5).Value = 10
$ws.Cells.Item(276, 6).Value = 100114014
$ws.Cells.Item(276, 7).Value = "Betarraga"
$ws.Cells.Item(276, 8).Value = "Sin especificar"
$ws.Cells.Item(276, 9).Value = "Primera"
$ws.Cells.Item(276, 10).Value = 500
$ws.Cells.Item(276, 11).Value = 1000
$ws.Cells.Item(276, 12).Value = 1000
$ws.Cells.Item(276, 13).Value = 1000
$ws.Cells.Item(276, 14).Value = "`$/paquete 5 unidades"
$ws.Cells.Item(276, 15).Value = "Región del Maule"
$ws.Cells.Item(276, 16).Value = 200
$ws.Cells.Item(276, 17).Value = 5
$ws.Cells.Item(276, 18).Value = "Hortaliza"
